# [Fix]: exclusion of 8 redundant metrics
#
# Removes the rows for 8 redundant metrics (MBRAE, UMBRAE, STDAPE, RMSPE,
# MRE, MRAE, MDRAE, GMRAE) from both worksheets, then renumbers the
# remaining "ID" column (column A) sequentially from 1.

$wb = $excel.ActiveWorkbook

$metricsToRemove = @("MBRAE", "UMBRAE", "STDAPE", "RMSPE", "MRE", "MRAE", "MDRAE", "GMRAE")

foreach ($ws in $wb.Worksheets) {
    $usedRows = $ws.UsedRange.Rows.Count

    # Walk bottom-up so deleting a row doesn't disturb the rows above it
    # that still need to be examined.
    for ($r = $usedRows; $r -ge 2; $r--) {
        $label = $ws.Cells.Item($r, 2).Value2
        if ($metricsToRemove -contains $label) {
            $ws.Rows.Item($r).Delete()
        }
    }

    # Renumber column A (the "ID" column) sequentially for the rows left.
    $usedRows = $ws.UsedRange.Rows.Count
    $id = 1
    for ($r = 2; $r -le $usedRows; $r++) {
        $ws.Cells.Item($r, 1).Value = $id
        $id = $id + 1
    }
}
